# Auto-generated Excel COM-interop script to apply the cryptos.xlsx symbol-list update
# Commit: "Updated symbol list on Tue Feb 14 23:54:24 UTC 2023 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextCell "D2" "296.32"
Set-TextCell "E2" "0.87%"
Set-TextCell "D3" "42.25"
Set-TextCell "E3" "3.73%"
Set-TextCell "D4" "5.029"
Set-TextCell "E4" "0.12%"
Set-TextCell "D5" "0.07582"
Set-TextCell "B6" "FTXToken"
Set-TextCell "C6" "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextCell "D6" "1.605"
Set-TextCell "E6" "3.22%"
Set-TextCell "B7" "MXToken"
Set-TextCell "C7" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextCell "D7" "0.9303"
Set-TextCell "E7" "0.64%"
Set-TextCell "B8" "BTSEToken"
Set-TextCell "C8" "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextCell "D8" "2.419"
Set-TextCell "E8" "2.98%"
Set-TextCell "B9" "LiechtensteinCryptoassetsExchange"
Set-TextCell "C9" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextCell "D9" "0.1212"
Set-TextCell "E9" "5.98%"
Set-TextCell "B10" "WazirX"
Set-TextCell "C10" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextCell "D10" "0.1841"
Set-TextCell "E10" "6.57%"
Set-TextCell "B11" "MandalaExchangeToken"
Set-TextCell "C11" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextCell "D11" "0.09002"
Set-TextCell "E11" "3.63%"
Set-TextCell "B12" "BitrueCoin"
Set-TextCell "C12" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextCell "D12" "0.04009"
Set-TextCell "E12" "-4.02%"
Set-TextCell "B13" "BitMartToken"
Set-TextCell "C13" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextCell "D13" "0.1053"
Set-TextCell "E13" "-0.13%"
Set-TextCell "B14" "BitForexToken"
Set-TextCell "C14" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextCell "D14" "0.001279"
Set-TextCell "E14" "1.15%"
Set-TextCell "B15" "TigerCash"
Set-TextCell "C15" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextCell "D15" "0.005792"
Set-TextCell "E15" "-1.89%"
Set-TextCell "B16" "LEO"
Set-TextCell "C16" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextCell "D16" "3.365"
Set-TextCell "E16" "-1.48%"
Set-TextCell "B17" "GateToken"
Set-TextCell "C17" "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextCell "D17" "4.393"
Set-TextCell "E17" "2.54%"
Set-TextCell "E18" "1.10%"
Set-TextCell "D19" "7.872"
Set-TextCell "E19" "2.38%"
Set-TextCell "D20" "0.1400"
Set-TextCell "E20" "1.67%"
Set-TextCell "E21" "4.25%"
Set-TextCell "D22" "0.04064"
Set-TextCell "E22" "5.16%"
Set-TextCell "D23" "0.001268"
Set-TextCell "E23" "0.78%"
Set-TextCell "D24" "0.003919"
Set-TextCell "E24" "1.02%"
Set-TextCell "D25" "0.0001231"
Set-TextCell "E25" "-3.76%"
Set-TextCell "E26" "0.15%"
Set-TextCell "D38" "0.02423"
Set-TextCell "E38" "3.67%"
Set-TextCell "D39" "0.05208"
Set-TextCell "E39" "3.96%"
Set-TextCell "D40" "0.006063"
Set-TextCell "E40" "10.85%"
Set-TextCell "D41" "0.007788"
Set-TextCell "E41" "1.15%"
Set-TextCell "D42" "0.1333"
Set-TextCell "D43" "0.007557"
Set-TextCell "E43" "2.96%"
Set-TextCell "D44" "0.007274"
Set-TextCell "E44" "-6.78%"
Set-TextCell "D45" "0.2964"
Set-TextCell "E45" "-6.30%"
Set-TextCell "D46" "0.00006781"
Set-TextCell "E46" "6.11%"
Set-TextCell "E47" "0.15%"
Set-TextCell "D48" "0.04503"
Set-TextCell "E48" "166.09%"
Set-TextCell "E50" "0.15%"
Set-TextCell "E51" "0.15%"
